$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the saved window position (best-effort; harmless if unsupported).
try { $excel.ActiveWindow.Left = 1160 } catch {}

# New rows 6-18 of to-do items (column A), plus one note in column B7.
$ws.Range("A6").Value  = "write out target criterion, expanding the terms so the signs become clear"
$ws.Range("A7").Value  = "mapping from Kalman gain to k_t --> compare estimation to evolution of Kalman gain"
$ws.Range("B7").Value  = "Gaetano Gaballo meeting, 11 June 2020, Notes 12, p 35"
$ws.Range("A8").Value  = "graph to understand the target criterion"
$ws.Range("A9").Value  = "two-period problem version of model w/ an intertemporal price people learn about, use it to show the mistake the CB makes when it assumes RE"
$ws.Range("A10").Value = "Gaetano's selling points:"
$ws.Range("A11").Value = "1. RE (a la Lucas) is great in long-run (don't fight RE!)"
$ws.Range("A12").Value = "2. but crisis has shown that things move quickly --> you'd need large changes in the model to match drifting long-run expectations"
$ws.Range("A13").Value = "3. so how costly is it to assume plain vanilla rational expectations in short run"
$ws.Range("A14").Value = "Recall that discretion=commitment w/o RE "
$ws.Range("A15").Value = "But how you depart from RE matters for how policy should deal with this fact. So I do 3 things"
$ws.Range("A16").Value = "1. estimate an adaptive learning model (small departures from RE can match empirical facts --> see how Marcet-Adam, Eusepi-Preston sell their models!)"
$ws.Range("A17").Value = "2. solve for optimal policy"
$ws.Range("A18").Value = "3. tell you what mistake CB makes when it instead assumes RE"

# Move the active selection to where editing left off.
$null = $ws.Range("A13").Select()
